# Documentación - Se actualiza plan de Calidad.
# Se añaden los documentos necesarios para completar la gestion de calidad:
# "Informe Final SQA", "Estándar de Programación"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renombrar los codigos de documentos de E1010-E1017 a E110-E117
$ws.Range("A11").Value = "E110"
$ws.Range("A12").Value = "E111"
$ws.Range("A13").Value = "E112"
$ws.Range("A14").Value = "E113"
$ws.Range("A15").Value = "E114"
$ws.Range("A16").Value = "E115"
$ws.Range("A17").Value = "E116"
$ws.Range("A18").Value = "E117"

# Se agrega el documento "Estándar de codificación" para E111
$ws.Range("B12").Value = "Estándar de codificación"

# "Plan de Pruebas" se traslada de la fila 13 a la fila 15 (E114)
$ws.Range("B13").ClearContents()
$ws.Range("B15").Value = "Plan de Pruebas"

# Actualizar la celda seleccionada tal como quedo en el archivo original
$ws.Range("B13").Select()
